# Edit commit: table style swap on 3 tables (slides 14-16) +
# slide-master theme colour scheme swap ("Integral"/Red Violet -> Office).
#
# The presentation ships two theme parts:
#   ppt/theme/theme1.xml -> used only by the Notes Master ("Office Theme")
#   ppt/theme/theme2.xml -> used by the (only) Slide Master ("Integral")
# The commit swaps the two themes' contents (color schemes) between the
# parts. The Notes Master's theme part is not reachable through the
# PowerPoint object model (Presentation.SlideMaster.Theme /
# Presentation.NotesMaster.Theme both resolve to the single theme that is
# actually wired to the Slide Master), so the only persisted, in-file
# change we can legitimately drive through COM is re-colouring the Slide
# Master's theme (ppt/theme/theme2.xml) to the palette the diff leaves
# there ("Office").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style: 3 tables (slides 14, 15, 16) move from the deck-local
#    "Table_0" style to the built-in style {D4014D65-F22D-4818-AA95-E53A00D17609}.
# ---------------------------------------------------------------------
$newTableStyle = "{D4014D65-F22D-4818-AA95-E53A00D17609}"
$tableSlides = @(14, 15, 16)

foreach ($slideIdx in $tableSlides) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Theme colour scheme: recolour the Slide Master's theme to the
#    "Office" palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
# ---------------------------------------------------------------------
function Set-ThemeHexRGB($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $scheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$officePalette = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($idx = 1; $idx -le 12; $idx++) {
    Set-ThemeHexRGB $themeColors $idx $officePalette[$idx - 1]
}
